{"js": "// Remove the trailing \"Requisitos\" heading paragraph and the\n// \"LOQ4240 - Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II (Requisito fraco)\" bullet\n// paragraph that follows the Bibliografia section at the end of the body.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (\n    text === \"Requisitos\" ||\n    text.indexOf(\"LOQ4240 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II  (Requisito fraco)\") !== -1\n  ) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Requisitos\" heading paragraph and the\n# \"LOQ4240 - Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II (Requisito fraco)\" bullet\n# paragraph that follows the Bibliografia section at the end of the body.\n\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"Requisitos\" -or $t.StartsWith(\"LOQ4240 -  Administra\u00e7\u00e3o e Organiza\u00e7\u00e3o II  (Requisito fraco)\")) {\n        $p.Range.Delete()\n    }\n}\n"}
